# local-vote-ppt.pptx: unify "Local Vote" wording with the app's "LocalVote"
# - Title slide headline: "LOCAL VOTE" -> "LOCALVOTE"
# - Title slide subtitle: "Presented by Team Local Vote" ->
#       "Presented by Team " + "LocalVote" (kept as a separate run)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the title ("Title 1") and subtitle ("Subtitle 2") placeholders on the
# first slide by name so the script is not strictly tied to shape ordering.
$title = $null
$subtitle = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "Title 1") { $title = $shp }
    if ($shp.Name -eq "Subtitle 2") { $subtitle = $shp }
}
if ($title -eq $null) { $title = $s.Shapes.Item(4) }
if ($subtitle -eq $null) { $subtitle = $s.Shapes.Item(5) }

# --- Title: "LOCAL VOTE" -> "LOCALVOTE" ---
$titleRange = $title.TextFrame.TextRange
$titleText = $titleRange.Text
$titleIdx = $titleText.IndexOf("LOCAL VOTE") + 1
if ($titleIdx -gt 0) {
    $titleRange.Characters($titleIdx, 10).Text = "LOCALVOTE"
}

# --- Subtitle: "...Local Vote" -> "...LocalVote" (becomes its own run) ---
$subRange = $subtitle.TextFrame.TextRange
$subText = $subRange.Text
$subIdx = $subText.IndexOf("Local Vote") + 1
if ($subIdx -gt 0) {
    $subRange.Characters($subIdx, 10).Text = "LocalVote"
}
